$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.68118155002594
$ws.Range("B1").Value = 1.870495915412903
$ws.Range("C1").Value = 1.937172293663025
$ws.Range("D1").Value = 2.51207160949707
$ws.Range("E1").Value = 3.432585954666138
